$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '89.749.35'
$ws.Range('E2').Value = '  +1.92%  '

$ws.Range('D3').Value = '3.188.65'
$ws.Range('E3').Value = '  -1.91%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.48'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.79%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '617.54'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.382'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.86%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.694'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.43%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.08%  '

$ws.Range('D10').Value = '3.191.88'
$ws.Range('E10').Value = '  -1.73%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.566'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.89%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.177'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.14%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.34'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.36%  '

$ws.Range('D15').Value = '3.779.54'
$ws.Range('E15').Value = '  -1.74%  '

$ws.Range('D16').Value = '89.627.45'
$ws.Range('E16').Value = '  +2.06%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '32.32'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.36%  '

$ws.Range('D18').Value = '3.207.81'
$ws.Range('E18').Value = '  -0.25%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000223'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +64.95%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.29'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +9.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.23'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.45%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '430.01'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.75%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.45'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.53%  '

$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.03'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.66%  '

$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.44'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.61%  '

$ws.Range('D27').Value = '3.366.28'
$ws.Range('E27').Value = '  -1.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '75.11'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.27%  '

$ws.Range('E29').Value = '  +0.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.22%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.152'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -14.91%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.98'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +30.94%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.32'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '530.84'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.88'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.76'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.25'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.94%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.07'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.43%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.31'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.04%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.125'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -8.42%  '

$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.89'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.67%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.366'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.85%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '149.93'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.84%  '

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.28'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.69%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '169.86'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.87%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.122'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -7.20%  '

$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.726'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.47%  '

$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.21'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.606'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.57%  '
